$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.608.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.656.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.04%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.69"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.70%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.656.33"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.16%  "

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "

$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.45%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.134.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.14%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.527.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.50%  "

$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.654.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.85%  "

$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  +0.34%  "

$ws.Range("E25").Value = "  +5.26%  "

$ws.Range("E26").Value = "  +6.26%  "

$ws.Range("E27").Value = "  +0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "549.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +19.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.01%  "

$ws.Range("E31").Value = "  -0.42%  "

$ws.Range("E32").Value = "  +14.23%  "

$ws.Range("E33").Value = "  +3.91%  "

$ws.Range("E34").Value = "  -0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "174.42"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.18%  "

$ws.Range("E36").Value = "  +9.42%  "

$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.95%  "

$ws.Range("E40").Value = "  +7.96%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "171.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.43%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("E43").Value = "  +2.14%  "

$ws.Range("E44").Value = "  +0.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "22.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.631"
$ws.Range("D46").Style = "Normal"

$ws.Range("E47").Value = "  +3.65%  "

$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("E49").Value = "  +1.94%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "

